$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.339.05"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.16"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.26"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4665"
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2737"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06294"
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.833.58"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07447"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.27"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.944"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "83.94"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6216"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.287.32"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.94"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007308"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.40"
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9995"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.926"
$ws.Range("E22").Value = "  -3.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.878"
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.196"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.95"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.86"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.878"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1026"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.077"
$ws.Range("E30").Value = "  -3.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.813"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04853"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.145"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7126"
$ws.Range("E34").Value = "  -2.17%  "
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.659"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8836"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "105.79"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.923"
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.550"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4035"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.101"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.64"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1201"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.644"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.29"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05513"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.356"
$ws.Range("E50").Value = "  -3.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3647"
$ws.Range("E51").Value = "  -1.49%  "
